$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "New Contract End Date" value for row 3 (column G) from 09/30/2022 to 09/30/2025.
# A leading apostrophe forces Excel to store the value as literal text instead
# of auto-converting the date-like string into a numeric date serial value,
# which keeps the existing cell style/number format untouched.
$ws.Range("G3").Value = "'09/30/2025"

# Update the selection to G3 as the active cell
$ws.Range("G3").Select()
